$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to be treated as text so numeric-looking
# strings (e.g. "230.96", "1.00") are not auto-converted to numbers.
$fmtRange = $ws.Range("D2:E51")
$fmtRange.NumberFormat = "@"

$ws.Range("D2").Value = '36.636.60'
$ws.Range("E2").Value = '  -2.26%  '
$ws.Range("D3").Value = '1.961.62'
$ws.Range("E3").Value = '  -3.19%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '230.96'
$ws.Range("E5").Value = '  -9.38%  '
$ws.Range("D6").Value = '0.597'
$ws.Range("E6").Value = '  -3.72%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '53.14'
$ws.Range("E8").Value = '  -6.99%  '
$ws.Range("D9").Value = '0.366'
$ws.Range("E9").Value = '  -4.94%  '
$ws.Range("D10").Value = '57.29'
$ws.Range("E10").Value = '  -0.03%  '
$ws.Range("D11").Value = '0.0736'
$ws.Range("E11").Value = '  -6.06%  '
$ws.Range("E12").Value = '  -4.31%  '
$ws.Range("D13").Value = '2.247.07'
$ws.Range("E13").Value = '  -3.41%  '
$ws.Range("D14").Value = '13.80'
$ws.Range("E14").Value = '  -5.20%  '
$ws.Range("D15").Value = '19.84'
$ws.Range("E15").Value = '  -6.55%  '
$ws.Range("D16").Value = '0.741'
$ws.Range("E16").Value = '  -9.32%  '
$ws.Range("B17").Value = 'Polkadot'
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D17").Value = '4.98'
$ws.Range("E17").Value = '  -7.03%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '1.971.28'
$ws.Range("E18").Value = '  -4.07%  '
$ws.Range("D19").Value = '36.548.92'
$ws.Range("E19").Value = '  -2.40%  '
$ws.Range("D20").Value = '67.18'
$ws.Range("E20").Value = '  -3.48%  '
$ws.Range("D21").Value = '0.0₃0792'
$ws.Range("E21").Value = '  -6.72%  '
$ws.Range("D22").Value = '4.98'
$ws.Range("E22").Value = '  -4.00%  '
$ws.Range("D23").Value = '220.27'
$ws.Range("E23").Value = '  -3.75%  '
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  +0.14%  '
$ws.Range("D25").Value = '2.33'
$ws.Range("E25").Value = '  -0.78%  '
$ws.Range("E26").Value = '  -11.09%  '
$ws.Range("D27").Value = '160.61'
$ws.Range("E27").Value = '  -2.19%  '
$ws.Range("D28").Value = '8.45'
$ws.Range("E28").Value = '  -6.62%  '
$ws.Range("D29").Value = '18.86'
$ws.Range("E29").Value = '  -5.21%  '
$ws.Range("D30").Value = '0.122'
$ws.Range("E30").Value = '  -6.05%  '
$ws.Range("D31").Value = '1.30'
$ws.Range("E31").Value = '  -5.11%  '
$ws.Range("D32").Value = '0.116'
$ws.Range("E32").Value = '  -3.69%  '
$ws.Range("D33").Value = '4.33'
$ws.Range("E33").Value = '  -8.12%  '
$ws.Range("D34").Value = '0.0602'
$ws.Range("E34").Value = '  -9.16%  '
$ws.Range("D35").Value = '4.17'
$ws.Range("E35").Value = '  -8.63%  '
$ws.Range("E36").Value = '  -7.19%  '
$ws.Range("E37").Value = '  -0.11%  '
$ws.Range("D39").Value = '3.17'
$ws.Range("E39").Value = '  -6.51%  '
$ws.Range("E40").Value = '  -3.12%  '
$ws.Range("E41").Value = '  +0.24%  '
$ws.Range("D42").Value = '1.406.50'
$ws.Range("E42").Value = '  +0.63%  '
$ws.Range("E43").Value = '  -7.00%  '
$ws.Range("D44").Value = '0.0869'
$ws.Range("E44").Value = '  -10.00%  '
$ws.Range("D45").Value = '1.08'
$ws.Range("E45").Value = '  -9.38%  '
$ws.Range("D46").Value = '86.49'
$ws.Range("E46").Value = '  -4.87%  '
$ws.Range("D47").Value = '0.982'
$ws.Range("E47").Value = '  -5.17%  '
$ws.Range("D48").Value = '14.74'
$ws.Range("E48").Value = '  -8.02%  '
$ws.Range("D49").Value = '2.86'
$ws.Range("E49").Value = '  -0.71%  '
$ws.Range("D50").Value = '6.66'
$ws.Range("E50").Value = '  -9.45%  '
$ws.Range("D51").Value = '2.139.65'
$ws.Range("E51").Value = '  -3.61%  '

# Restore the original (default) style so no residual text-format
# styling differs from the source workbook.
$fmtRange.Style = "Normal"

